$d = $word.ActiveDocument

$replacements = @(
    @("370÷8=", "850÷2="),
    @("464÷7=", "984÷9="),
    @("489÷3=", "392÷5="),
    @("914÷6=", "931÷8="),
    @("865÷7=", "648÷6="),
    @("274÷4=", "792÷9="),
    @("685÷4=", "915÷7="),
    @("384÷6=", "271÷4="),
    @("207÷2=", "161÷9="),
    @("816÷7=", "869÷6="),
    @("605÷5=", "788÷5="),
    @("853÷6=", "825÷3="),
    @("217÷9=", "362÷3="),
    @("662÷8=", "675÷4="),
    @("828÷7=", "573÷4="),
    @("483÷4=", "310÷4="),
    @("107÷3=", "660÷9="),
    @("455÷5=", "125÷7="),
    @("968÷2=", "623÷6="),
    @("528÷8=", "964÷3="),
    @("394÷8=", "265÷4="),
    @("230÷7=", "433÷5="),
    @("804÷2=", "472÷5="),
    @("154÷7=", "527÷7="),
    @("402÷8=", "793÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
